$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data-source filename referenced in column A (rows 2-5)
$ws.Range("A2").Value = "patients_2_specification.xlsx"
$ws.Range("A3").Value = "patients_2_specification.xlsx"
$ws.Range("A4").Value = "patients_2_specification.xlsx"
$ws.Range("A5").Value = "patients_2_specification.xlsx"

# Update the selected/active cell to A2
[void]$ws.Range("A2").Select()

# Touch page setup (portrait) to mirror the resave side effect captured in the diff
$ws.PageSetup.Orientation = 1
